$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cell ranges that become blank (previously all zero) - representing rows with
# concatenated/blank quarterly columns that have no data for this balance sheet.
$blankRanges = "E57:S57", "E58:S58", "C64:E64", "E71:S71", "E72:S72", "E73:S73", "E77:S77", "E78:S78"
foreach ($rngAddr in $blankRanges) {
    $rng = $ws.Range($rngAddr)
    $rng.Value = ""
    $rng.Style = "Normal"
}

# Update recomputed values for column E (and a couple of others) reflecting the
# newly concatenated balance sheet figures.
$ws.Range("E59").Value = 390490.88
$ws.Range("M59").Value = 448109.024
$ws.Range("Q59").Value = 932511.1040000001
$ws.Range("E60").Value = -283268.992
$ws.Range("I60").Value = -282228.032
$ws.Range("E61").Value = 107221.984
$ws.Range("I61").Value = 77098
$ws.Range("M61").Value = 97650.984
$ws.Range("E62").Value = -5042
$ws.Range("E63").Value = -18440
$ws.Range("E66").Value = -29771
$ws.Range("E68").Value = -10044
$ws.Range("E69").Value = 8490
$ws.Range("E70").Value = -18534
$ws.Range("E74").Value = 43925.008
$ws.Range("E75").Value = -4642
$ws.Range("E76").Value = -1920
$ws.Range("E80").Value = 29954
